$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells remain text (avoid numeric auto-conversion)
$dCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "26.403.08"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.792.62"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "307.27"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").Value = "0.4560"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").Value = "0.3621"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "0.8785"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "0.07846"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.275"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.731.91"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("D15").Value = "6.321"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "84.94"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "0.000008557"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "26.436.77"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "14.27"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "4.984"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "10.47"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "1.977.16"
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "17.91"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "2.036"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "111.96"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "4.857"
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "0.08655"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "3.073"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "4.445"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("D35").Value = "2.666"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("D37").Value = "1.005"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "1.076"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "0.01942"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").Value = "0.05113"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "2.872"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "0.5228"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("D43").Value = "6.890"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -4.58%  "
$ws.Range("D45").Value = "8.019"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").Value = "0.4706"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "9.865"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").Value = "99.91"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").Value = "1.585"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "0.05963"
$ws.Range("E51").Value = "  -2.14%  "
